# Thesis roadmap figure touch-up:
#  - nudge a large number of boxes/connectors horizontally (EMU-precise)
#  - restack three shapes (id 18, 627, 629) to the front of the z-order

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# New Left (in points) for each shape, keyed by shape Id. PowerPoint's
# Shape.Left is a single-precision (float32) property, and this host
# truncates (rather than rounds) points->EMU on write-back, so the values
# below are nudged by a sub-ulp amount versus the naive EMU/12700 so that
# they still resolve to the exact target EMU once stored.
$newLeft = @{
    2   = 98.92330708661417
    3   = 98.92330708661417
    4   = 98.92330708661417
    5   = 98.92330708661417
    6   = 22.92527559055118
    7   = 124.16937007874016
    9   = 98.92330708661417
    10  = 98.92330708661417
    11  = 103.17527559055118
    12  = 103.17527559055118
    13  = 103.17527559055118
    14  = 103.17527559055118
    15  = 103.17527559055118
    16  = 103.17527559055118
    17  = 126.95417412834647
    19  = 26.494094488188978
    28  = 158.8792190984252
    29  = 174.04142002283464
    30  = 174.04142002283464
    31  = 174.04142002283464
    32  = 74.34299212598425
    33  = 173.8644094488189
    35  = 173.39275590551182
    36  = 174.04142002283464
    41  = 174.04142002283464
    46  = 48.518740157480316
    47  = 51.11338622677165
    48  = 51.11338622677165
    49  = 51.11338622677165
    50  = 26.42220472440945
    51  = 125.09669291338582
    53  = 51.11338622677165
    54  = 51.11338622677165
    58  = 193.12070866141732
    59  = 193.20063022125984
    60  = 259.3652801905512
    591 = 38.927795275590555
    594 = 26.494094488188978
    595 = 26.494094488188978
    596 = 26.494094488188978
    597 = 26.494094488188978
    598 = 37.027166454330704
    600 = 25.12488278976378
    619 = 126.95417412834647
    620 = 126.95417412834647
    622 = 126.95417412834647
    623 = 126.95417412834647
    624 = 136.99409448818898
    625 = 126.95417412834647
    626 = 140.6088188976378
    631 = 226.76551181102363
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($newLeft.ContainsKey($sh.Id)) {
        $sh.Left = $newLeft[$sh.Id]
    }
}

# Helper to find a shape on the slide by its Id.
function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

# Bring shapes 18, 627, 629 to the front of the z-order (in this order),
# matching the author's restack of the "Research Question 3" / "Main output(s)"
# boxes and the empty "Rectangle 17" placeholder.
(Get-ShapeById $s 18).ZOrder(0)
(Get-ShapeById $s 627).ZOrder(0)
(Get-ShapeById $s 629).ZOrder(0)
